# [#PAB-232] remove partially left-over 'sheet_names' concept in ingest endpoint
#
# - Add a new worksheet ("Another Sheet") after the existing "TestSheet",
#   populate it with a header row ("Field 1"/"Field 2") and fill the
#   remaining cells with "Data".
# - Make the new sheet the active/selected sheet, with B4 as the active cell.
# - Update the selection on the original "TestSheet" to A1:C1 (and it is no
#   longer the tab-selected sheet, since focus moved to the new sheet).

$wb = $excel.ActiveWorkbook

# Original (and so far, only) worksheet.
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet directly after TestSheet, so the sheet order becomes
# TestSheet, Another Sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Another Sheet"

# Header row.
$ws2.Range("A1").Value = "Field 1"
$ws2.Range("B1").Value = "Field 2"

# Remaining data cells.
$ws2.Range("A2:B4").Value = "Data"

# Update the selection shown on TestSheet.
[void]$ws1.Range("A1:C1").Select()

# Make "Another Sheet" the active sheet/tab, with B4 selected - this is the
# last sheet/selection touched, so it ends up as the active one on save.
$ws2.Activate()
[void]$ws2.Range("B4").Select()
